$d = $word.ActiveDocument

# 1. Bold the two heading paragraphs
$d.Paragraphs.Item(1).Range.Bold = 1
$d.Paragraphs.Item(3).Range.Bold = 1

# 2. Merge/clean runs that had proofErr spell-check markers
$d.Content.Find.Execute("Freescale ColdFire (32-bit) và S08 (8-bit)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Freescale ColdFire (32-bit) và S08 (8-bit)", 2) | Out-Null

$d.Content.Find.Execute("- Renesas Electronics: RL78 16-bit MCU; RX 32-bit MCU; SuperH; V850 32-", $true, $false, $false, $false, $false,
                         $true, 1, $false, "- Renesas Electronics: RL78 16-bit MCU; RX 32-bit MCU; SuperH; V850 32-", 2) | Out-Null

$d.Content.Find.Execute("- PSoC (Programmable System-on-Chip)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "- PSoC (Programmable System-on-Chip)", 2) | Out-Null

# 3. Replace the "Stellaris" paragraph (merging its runs/removing proofErr + stray bookmark)
#    and inject all of the new RS485 section content right after it in one shot.
$stellarisPara = $d.Paragraphs.Item(13)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Stellaris (32-bit)</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>RS485</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Introduction</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  Nowadays, there are various communication protocols can be used for the thesis, namely I2C, ISP, RS232, RS485, Bluetooth or Wi-Fi. Each protocol is designed to be suitable for specified purpose with different advantages or disadvantages, which means a perfect protocol does not exist. </w:t></w:r><w:r><w:t xml:space="preserve">When making a decision to choose suitable protocols for the thesis, the author had to think about the trade-off between the stabilization </w:t></w:r><w:r><w:t>and the speed of the communication protocol.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">RS485 is chosen as the main way for components in the system to communicate with </w:t></w:r><w:r><w:t>each other</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>RS485</w:t></w:r><w:r><w:t xml:space="preserve"> is defined in 1983 not as a protocol but an </w:t></w:r><w:r><w:t xml:space="preserve">electrical </w:t></w:r><w:r><w:t xml:space="preserve">interface standard and only specifies the </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>drivers and receivers’ characteristics.</w:t></w:r><w:r><w:t xml:space="preserve"> It is developed in order to make data rate and transmitting dis</w:t></w:r><w:r><w:t>tance are inversely proportional. For instance, the data transmitting speed can reach 10 Mbps within distance of 16 meters or if the distance is extended to 1220 meters, the data rate is lower to 100 kbps.</w:t></w:r><w:r><w:t xml:space="preserve"> The advantage of RS485 over RS232, which is developed in 1960, is multiple nodes can be parallel connected to a bus</w:t></w:r><w:r><w:t>. Additionally, the network can be extended in length and number of nodes easily by using simple connector.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Besides, Wi-Fi, Bluetooth and </w:t></w:r><w:r><w:t>GSM</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> are</w:t></w:r><w:r><w:t xml:space="preserve"> also implemented in the thesis in order to take the advantages in different circumstances. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>RS485 specification</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Overview</w:t></w:r></w:p><w:p><w:r><w:t>\table</w:t></w:r></w:p><w:p><w:r><w:t>Table shows the highlight specifications of RS485.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$stellarisPara.Range.InsertXML($xml) | Out-Null

Write-Output "done"
